# Apply cryptos.xlsx price/volume/coin updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '92.440.60'
$ws.Range("E2").Value = '  +2.35%  '
# Row 3
$ws.Range("D3").Value = '3.153.04'
$ws.Range("E3").Value = '  +2.20%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.64%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.99'
$ws.Range("E5").Value = '  +4.23%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.60'
$ws.Range("E6").Value = '  +0.15%  '
# Row 7
$ws.Range("E7").Value = '  +8.43%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.378'
$ws.Range("E8").Value = '  +3.83%  '
# Row 9
$ws.Range("E9").Value = '  -0.16%  '
# Row 10
$ws.Range("D10").Value = '3.148.18'
$ws.Range("E10").Value = '  +2.09%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.766'
$ws.Range("E11").Value = '  +5.80%  '
# Row 12
$ws.Range("E12").Value = '  +4.54%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value = '  -0.53%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.64'
$ws.Range("E14").Value = '  -2.27%  '
# Row 15
$ws.Range("E15").Value = '  +1.58%  '
# Row 16
$ws.Range("D16").Value = '92.048.75'
$ws.Range("E16").Value = '  +2.37%  '
# Row 17
$ws.Range("D17").Value = '3.721.70'
$ws.Range("E17").Value = '  +1.57%  '
# Row 18
$ws.Range("D18").Value = '3.109.71'
$ws.Range("E18").Value = '  +0.33%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.77'
$ws.Range("E19").Value = '  -1.32%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.15'
$ws.Range("E20").Value = '  +8.87%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000214'
$ws.Range("E21").Value = '  -0.37%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.87'
$ws.Range("E22").Value = '  +5.93%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '447.53'
$ws.Range("E23").Value = '  +2.63%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.25'
$ws.Range("E24").Value = '  +4.11%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.73'
$ws.Range("E25").Value = '  +0.37%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.75'
$ws.Range("E26").Value = '  +9.01%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.05'
$ws.Range("E27").Value = '  -0.79%  '
# Row 28
$ws.Range("D28").Value = '3.246.34'
$ws.Range("E28").Value = '  -0.22%  '
# Row 29
$ws.Range("E29").Value = '  +0.00%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.250'
$ws.Range("E30").Value = '  +27.84%  '
# Row 31
$ws.Range("E31").Value = '  +10.86%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.125'
$ws.Range("E32").Value = '  +42.85%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.37'
$ws.Range("E33").Value = '  -0.08%  '
# Row 34
$ws.Range("E34").Value = '  +11.33%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.06'
$ws.Range("E35").Value = '  +12.78%  '
# Row 36
$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.937'
$ws.Range("E36").Value = '  -5.97%  '
# Row 37
$ws.Range("B37").Value = 'MantraDAO'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.35'
$ws.Range("E37").Value = '  +25.88%  '
# Row 38
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.76'
$ws.Range("E38").Value = '  +4.19%  '
# Row 39
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '502.72'
$ws.Range("E39").Value = '  -0.50%  '
# Row 40
$ws.Range("B40").Value = 'PancakeSwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.94'
$ws.Range("E40").Value = '  +2.17%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.59'
$ws.Range("E41").Value = '  -7.11%  '
# Row 42
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.32'
$ws.Range("E42").Value = '  +2.55%  '
# Row 43
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.428'
$ws.Range("E43").Value = '  +5.07%  '
# Row 44
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.27'
$ws.Range("E44").Value = '  +0.44%  '
# Row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.01%  '
# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.95'
$ws.Range("E46").Value = '  +2.84%  '
# Row 47
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.708'
$ws.Range("E47").Value = '  +3.33%  '
# Row 48
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '154.43'
$ws.Range("E48").Value = '  +1.98%  '
# Row 49
$ws.Range("E49").Value = '  +1.93%  '
# Row 50
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.51'
$ws.Range("E50").Value = '  -0.89%  '
# Row 51
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.46'
$ws.Range("E51").Value = '  +2.05%  '
